$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (mirrors the existing row layout: A..T)
$rows = @(
    @{
        A = 9
        B = "Vega Central Mapocho de Santiago"
        C = "Metropolitana"
        D = 44595
        E = 13
        F = "Fruta"
        G = 100103
        H = "Frutos de hueso (carozo)"
        I = 100103002
        J = "Ciruela"
        K = "Black Amber"
        L = "Especial"
        M = 310
        N = 10500
        O = 10500
        P = 10500
        Q = "`$/caja 15 kilos granel"
        R = "Región de O'Higgins"
        S = 700
        T = 15
    },
    @{
        A = 9
        B = "Vega Central Mapocho de Santiago"
        C = "Metropolitana"
        D = 44595
        E = 13
        F = "Fruta"
        G = 100103
        H = "Frutos de hueso (carozo)"
        I = 100103002
        J = "Ciruela"
        K = "Black Amber"
        L = "Primera"
        M = 480
        N = 7500
        O = 7500
        P = 7500
        Q = "`$/caja 15 kilos granel"
        R = "Región de O'Higgins"
        S = 500
        T = 15
    }
)

$startRow = 86

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data.A
    $ws.Cells.Item($r, 2).Value = $data.B
    $ws.Cells.Item($r, 3).Value = $data.C

    # Column D carries the same date-time number format as the rows above it
    $ws.Cells.Item($r, 4).Value = $data.D
    $ws.Cells.Item($r, 4).NumberFormat = $ws.Cells.Item($r - 1, 4).NumberFormat

    $ws.Cells.Item($r, 5).Value = $data.E
    $ws.Cells.Item($r, 6).Value = $data.F
    $ws.Cells.Item($r, 7).Value = $data.G
    $ws.Cells.Item($r, 8).Value = $data.H
    $ws.Cells.Item($r, 9).Value = $data.I
    $ws.Cells.Item($r, 10).Value = $data.J
    $ws.Cells.Item($r, 11).Value = $data.K
    $ws.Cells.Item($r, 12).Value = $data.L
    $ws.Cells.Item($r, 13).Value = $data.M
    $ws.Cells.Item($r, 14).Value = $data.N
    $ws.Cells.Item($r, 15).Value = $data.O
    $ws.Cells.Item($r, 16).Value = $data.P
    $ws.Cells.Item($r, 17).Value = $data.Q
    $ws.Cells.Item($r, 18).Value = $data.R
    $ws.Cells.Item($r, 19).Value = $data.S
    $ws.Cells.Item($r, 20).Value = $data.T
}
